# 505.1 Double Columns (::) Operator - fix title typo "Opeator" -> "Operator"
# and merge the two runs back into a single run (matching the author's
# manual correction of the misspelled word).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# The title text is currently split into two runs:
#   "505.1 Double Columns (::) " + "Opeator" (marked err="1" - spellcheck flag)
# Remove the misspelled second run's characters entirely, then append the
# correctly spelled word onto the end of the (now sole) run so PowerPoint
# folds everything into one run using the first run's formatting.
$badWord = $tr.Characters(27, 7)
$badWord.Text = ""

[void]$tr.InsertAfter("Operator")

# Re-apply the full run's text onto itself so the inserted text merges into
# a single run with the original formatting (size 4000, bold, yellow fill)
# instead of staying a separate run.
$full = $tr.Characters(1, $tr.Length)
$full.Text = "505.1 Double Columns (::) Operator"
